# Actualización automática del tracker
# Marca resultados "Fallo" (profit -1) para las filas de eventos ya resueltos
# que aún no tenían resultado registrado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(25, 26, 40, 43, 47, 48)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Fallo"
    $ws.Range("H$r").Value = -1
}
